$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.484760527490852
$ws.Range("C2").Value = 2.974931511405659
$ws.Range("D2").Value = 73.29250480412391
$ws.Range("E2").Value = 0.001200297787319263
$ws.Range("F2").Value = 0.3070136037349682
$ws.Range("G2").Value = -1.576088003932502
$ws.Range("H2").Value = -0.5252816652720125
$ws.Range("I2").Value = 1.258298382116306
$ws.Range("J2").Value = 4.353833897097513
$ws.Range("K2").Value = 933
$ws.Range("L2").Value = -27.01326210292276
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 4.355597361710274
$ws.Range("O2").Value = 5.022264028376941

$ws.Range("B3").Value = 215.8460488256372
$ws.Range("C3").Value = 54.41977192452163
$ws.Range("D3").Value = 0.06296368512311659
$ws.Range("E3").Value = -154.2868618104845
$ws.Range("F3").Value = -0.329547730428837
$ws.Range("G3").Value = 0.07897087730742713
$ws.Range("H3").Value = 0.8249624745230482
$ws.Range("I3").Value = -0.2167969128771863
$ws.Range("J3").Value = 4.353822396771477
$ws.Range("K3").Value = 365
$ws.Range("L3").Value = -65.00020152445437
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 4.35562456493966
$ws.Range("O3").Value = 5.022291231606327

$ws.Range("B4").Value = 0.0001375871783498155
$ws.Range("C4").Value = 111.3254406937106
$ws.Range("D4").Value = 26.02470625388843
$ws.Range("E4").Value = 97.09511874251959
$ws.Range("F4").Value = 1.749804386321745
$ws.Range("G4").Value = -0.331850500117137
$ws.Range("H4").Value = -1.765056258813615
$ws.Range("I4").Value = 0.08660936900193228
$ws.Range("J4").Value = 4.35382900859284
$ws.Range("K4").Value = 104
$ws.Range("L4").Value = -162.6976376054434
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 4.355700516939208
$ws.Range("O4").Value = 5.022367183605875

$ws.Range("B5").Value = 62.00401219989043
$ws.Range("C5").Value = -0.0686988724152805
$ws.Range("D5").Value = 16.81072839481216
$ws.Range("E5").Value = 2.596492926043097
$ws.Range("F5").Value = -0.742261392532295
$ws.Range("G5").Value = 0.6403693976631497
$ws.Range("H5").Value = -0.2935299578201138
$ws.Range("I5").Value = 0.4431989330853714
$ws.Range("J5").Value = 4.353829226599043
$ws.Range("K5").Value = 949
$ws.Range("L5").Value = -18.48999883258931
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 4.3557377345204
$ws.Range("O5").Value = 5.022404401187067

$ws.Range("B6").Value = 0.001027507922224007
$ws.Range("C6").Value = 1.031840058180127
$ws.Range("D6").Value = 28.58718385622127
$ws.Range("E6").Value = 87.03344124046752
$ws.Range("F6").Value = 1.447963830809482
$ws.Range("G6").Value = -1.229301218265868
$ws.Range("H6").Value = 0.1471364779517068
$ws.Range("I6").Value = -0.5120584981971066
$ws.Range("J6").Value = 4.353740071341456
$ws.Range("K6").Value = 524
$ws.Range("L6").Value = -58.77571923202518
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 4.355760881489715
$ws.Range("O6").Value = 5.022427548156382

$ws.Range("B7").Value = 46.44902483438472
$ws.Range("C7").Value = 37.09655258550981
$ws.Range("D7").Value = 0.1384490193450018
$ws.Range("E7").Value = 3.851220908770927
$ws.Range("F7").Value = -0.9018521530548023
$ws.Range("G7").Value = -0.4844644207146638
$ws.Range("H7").Value = 0.6954831642038024
$ws.Range("I7").Value = 0.3215955571925955
$ws.Range("J7").Value = 4.353801292017199
$ws.Range("K7").Value = 585
$ws.Range("L7").Value = -18.48980883694882
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 4.355842516964856
$ws.Range("O7").Value = 5.022509183631523

$ws.Range("B8").Value = 74.87480888822702
$ws.Range("C8").Value = -0.421800843465856
$ws.Range("D8").Value = 95.24470402266655
$ws.Range("E8").Value = 0.1115279975677855
$ws.Range("F8").Value = 0.08076027322490242
$ws.Range("G8").Value = 0.645623834798482
$ws.Range("H8").Value = -0.5332206045361556
$ws.Range("I8").Value = 0.8946585142796977
$ws.Range("J8").Value = 4.35371576955643
$ws.Range("K8").Value = 721
$ws.Range("L8").Value = -108.837569616994
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 4.355846699701061
$ws.Range("O8").Value = 5.022513366367728

$ws.Range("B9").Value = -25.30642917605476
$ws.Range("C9").Value = 70.56436989947386
$ws.Range("D9").Value = 3.981789824980415
$ws.Range("E9").Value = 68.18925575230816
$ws.Range("F9").Value = -0.188233701045408
$ws.Range("G9").Value = -1.040866130936773
$ws.Range("H9").Value = 0.3883552408162663
$ws.Range("I9").Value = -0.1953932198606885
$ws.Range("J9").Value = 4.353839073164195
$ws.Range("K9").Value = 376
$ws.Range("L9").Value = -34.92325937647293
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 4.355936935495713
$ws.Range("O9").Value = 5.02260360216238

$ws.Range("B10").Value = 119.7326409213468
$ws.Range("C10").Value = 81.94047024097438
$ws.Range("D10").Value = 85.03416273998866
$ws.Range("E10").Value = 0.0007087102803752547
$ws.Range("F10").Value = 0.06907366095363976
$ws.Range("G10").Value = -0.07227372586296976
$ws.Range("H10").Value = -0.5927917729632113
$ws.Range("I10").Value = 1.55055399719802
$ws.Range("J10").Value = 4.353691795239314
$ws.Range("K10").Value = 451
$ws.Range("L10").Value = -223.2023377832309
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 4.355944973611741
$ws.Range("O10").Value = 5.022611640278408

$ws.Range("B11").Value = 0.5753863727734907
$ws.Range("C11").Value = 0.1205574064483508
$ws.Range("D11").Value = 85.2498340834572
$ws.Range("E11").Value = 8.632410031552443
$ws.Range("F11").Value = 0.2522334519041363
$ws.Range("G11").Value = 0.7961275668097332
$ws.Range("H11").Value = -0.7589242339494804
$ws.Range("I11").Value = 0.1628426857156886
$ws.Range("J11").Value = 4.353735884756489
$ws.Range("K11").Value = 534
$ws.Range("L11").Value = -20.85345408414249
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 4.355987731520631
$ws.Range("O11").Value = 5.022784483094951
